$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the rows for department codes that are no longer wanted:
# BA (row 2), CS (row 4), ENG (row 5), GD (row 7), HS (row 8), TI (row 14)
# Deleted from the bottom up so earlier row numbers stay valid.
$ws.Rows(14).Delete()
$ws.Rows(8).Delete()
$ws.Rows(7).Delete()
$ws.Rows(5).Delete()
$ws.Rows(4).Delete()
$ws.Rows(2).Delete()

# Re-order the remaining department codes, moving SS to the top of the list
$codes = @("SS", "CF", "FSOFT", "ITS", "JPN", "Math", "SE")
for ($i = 0; $i -lt $codes.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $codes[$i]
    $ws.Cells.Item($row, 2).Value = $codes[$i]
}

# Update the selected cell to match the last edited location
$ws.Range("B8").Select()
